$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "Switzerland"
$ws.Cells.Item(2,2).Value = 1889
$ws.Cells.Item(2,3).Value = 928
$ws.Cells.Item(2,4).Value = 1148
$ws.Cells.Item(2,5).Value = -1667
$ws.Cells.Item(2,6).Value = 3909
$ws.Cells.Item(2,7).Value = 1345
$ws.Cells.Item(2,8).Value = -1548
$ws.Cells.Item(2,9).Value = 4005

$ws.Cells.Item(3,1).Value = "Switzerland"
$ws.Cells.Item(3,2).Value = 1890
$ws.Cells.Item(3,3).Value = 3408
$ws.Cells.Item(3,4).Value = 3356
$ws.Cells.Item(3,5).Value = 541
$ws.Cells.Item(3,6).Value = 6078
$ws.Cells.Item(3,7).Value = 3451
$ws.Cells.Item(3,8).Value = 711
$ws.Cells.Item(3,9).Value = 6123

$ws.Cells.Item(4,1).Value = "Switzerland"
$ws.Cells.Item(4,2).Value = 1891
$ws.Cells.Item(4,3).Value = 2024
$ws.Cells.Item(4,4).Value = 1593
$ws.Cells.Item(4,5).Value = -1471
$ws.Cells.Item(4,6).Value = 4492
$ws.Cells.Item(4,7).Value = 1955
$ws.Cells.Item(4,8).Value = -1328
$ws.Cells.Item(4,9).Value = 4976

$ws.Cells.Item(5,1).Value = "Switzerland"
$ws.Cells.Item(5,2).Value = 1917
$ws.Cells.Item(5,3).Value = 3621
$ws.Cells.Item(5,4).Value = 3833
$ws.Cells.Item(5,5).Value = 1138
$ws.Cells.Item(5,6).Value = 6502
$ws.Cells.Item(5,7).Value = 3832
$ws.Cells.Item(5,8).Value = 1212
$ws.Cells.Item(5,9).Value = 6476

$ws.Cells.Item(6,1).Value = "Switzerland"
$ws.Cells.Item(6,2).Value = 1918
$ws.Cells.Item(6,3).Value = 24562
$ws.Cells.Item(6,4).Value = 24692
$ws.Cells.Item(6,5).Value = 21817
$ws.Cells.Item(6,6).Value = 27411
$ws.Cells.Item(6,7).Value = 24779
$ws.Cells.Item(6,8).Value = 21958
$ws.Cells.Item(6,9).Value = 27548

$ws.Cells.Item(7,1).Value = "Switzerland"
$ws.Cells.Item(7,2).Value = 1919
$ws.Cells.Item(7,3).Value = 3763
$ws.Cells.Item(7,4).Value = 3654
$ws.Cells.Item(7,5).Value = 352
$ws.Cells.Item(7,6).Value = 6730
$ws.Cells.Item(7,7).Value = 3750
$ws.Cells.Item(7,8).Value = 437
$ws.Cells.Item(7,9).Value = 6918

$ws.Cells.Item(8,1).Value = "Switzerland"
$ws.Cells.Item(8,2).Value = 1956
$ws.Cells.Item(8,3).Value = 1556
$ws.Cells.Item(8,4).Value = 1470
$ws.Cells.Item(8,5).Value = -893
$ws.Cells.Item(8,6).Value = 3756
$ws.Cells.Item(8,7).Value = 1513
$ws.Cells.Item(8,8).Value = -778
$ws.Cells.Item(8,9).Value = 3744

$ws.Cells.Item(9,1).Value = "Switzerland"
$ws.Cells.Item(9,2).Value = 1957
$ws.Cells.Item(9,3).Value = -1215
$ws.Cells.Item(9,4).Value = -1046
$ws.Cells.Item(9,5).Value = -3859
$ws.Cells.Item(9,6).Value = 1658
$ws.Cells.Item(9,7).Value = -961
$ws.Cells.Item(9,8).Value = -3812
$ws.Cells.Item(9,9).Value = 1632

$ws.Cells.Item(10,1).Value = "Switzerland"
$ws.Cells.Item(10,2).Value = 1958
$ws.Cells.Item(10,3).Value = -3376
$ws.Cells.Item(10,4).Value = -3234
$ws.Cells.Item(10,5).Value = -6396
$ws.Cells.Item(10,6).Value = -148
$ws.Cells.Item(10,7).Value = -3603
$ws.Cells.Item(10,8).Value = -6833
$ws.Cells.Item(10,9).Value = -595

$ws.Cells.Item(11,1).Value = "Switzerland"
$ws.Cells.Item(11,2).Value = 2019
$ws.Cells.Item(11,3).Value = 233
$ws.Cells.Item(11,4).Value = 315
$ws.Cells.Item(11,5).Value = -2991
$ws.Cells.Item(11,6).Value = 3505
$ws.Cells.Item(11,7).Value = 202
$ws.Cells.Item(11,8).Value = -2994
$ws.Cells.Item(11,9).Value = 3370

$ws.Cells.Item(12,1).Value = "Switzerland"
$ws.Cells.Item(12,2).Value = 2020
$ws.Cells.Item(12,3).Value = 8123
$ws.Cells.Item(12,4).Value = 7959
$ws.Cells.Item(12,5).Value = 4808
$ws.Cells.Item(12,6).Value = 11009
$ws.Cells.Item(12,7).Value = 9478
$ws.Cells.Item(12,8).Value = 6380
$ws.Cells.Item(12,9).Value = 12482

$ws.Cells.Item(13,1).Value = "Sweden"
$ws.Cells.Item(13,2).Value = 1889
$ws.Cells.Item(13,3).Value = 2473
$ws.Cells.Item(13,4).Value = 2208
$ws.Cells.Item(13,5).Value = -1577
$ws.Cells.Item(13,6).Value = 5915
$ws.Cells.Item(13,7).Value = 2448
$ws.Cells.Item(13,8).Value = -1396
$ws.Cells.Item(13,9).Value = 5970

$ws.Cells.Item(14,1).Value = "Sweden"
$ws.Cells.Item(14,2).Value = 1890
$ws.Cells.Item(14,3).Value = 8747
$ws.Cells.Item(14,4).Value = 8514
$ws.Cells.Item(14,5).Value = 4555
$ws.Cells.Item(14,6).Value = 12370
$ws.Cells.Item(14,7).Value = 8782
$ws.Cells.Item(14,8).Value = 4830
$ws.Cells.Item(14,9).Value = 12435

$ws.Cells.Item(15,1).Value = "Sweden"
$ws.Cells.Item(15,2).Value = 1891
$ws.Cells.Item(15,3).Value = 6268
$ws.Cells.Item(15,4).Value = 5405
$ws.Cells.Item(15,5).Value = 822
$ws.Cells.Item(15,6).Value = 9767
$ws.Cells.Item(15,7).Value = 5374
$ws.Cells.Item(15,8).Value = 976
$ws.Cells.Item(15,9).Value = 9670

$ws.Cells.Item(16,1).Value = "Sweden"
$ws.Cells.Item(16,2).Value = 1917
$ws.Cells.Item(16,3).Value = -3025
$ws.Cells.Item(16,4).Value = -2935
$ws.Cells.Item(16,5).Value = -8320
$ws.Cells.Item(16,6).Value = 2093
$ws.Cells.Item(16,7).Value = -2695
$ws.Cells.Item(16,8).Value = -7899
$ws.Cells.Item(16,9).Value = 2294

$ws.Cells.Item(17,1).Value = "Sweden"
$ws.Cells.Item(17,2).Value = 1918
$ws.Cells.Item(17,3).Value = 25949
$ws.Cells.Item(17,4).Value = 26037
$ws.Cells.Item(17,5).Value = 20771
$ws.Cells.Item(17,6).Value = 31224
$ws.Cells.Item(17,7).Value = 26022
$ws.Cells.Item(17,8).Value = 20898
$ws.Cells.Item(17,9).Value = 30895

$ws.Cells.Item(18,1).Value = "Sweden"
$ws.Cells.Item(18,2).Value = 1919
$ws.Cells.Item(18,3).Value = 8440
$ws.Cells.Item(18,4).Value = 7293
$ws.Cells.Item(18,5).Value = 1320
$ws.Cells.Item(18,6).Value = 13151
$ws.Cells.Item(18,7).Value = 7296
$ws.Cells.Item(18,8).Value = 1279
$ws.Cells.Item(18,9).Value = 13184

$ws.Cells.Item(19,1).Value = "Sweden"
$ws.Cells.Item(19,2).Value = 1956
$ws.Cells.Item(19,3).Value = 1574
$ws.Cells.Item(19,4).Value = 1397
$ws.Cells.Item(19,5).Value = -2176
$ws.Cells.Item(19,6).Value = 4740
$ws.Cells.Item(19,7).Value = 1050
$ws.Cells.Item(19,8).Value = -2609
$ws.Cells.Item(19,9).Value = 4433

$ws.Cells.Item(20,1).Value = "Sweden"
$ws.Cells.Item(20,2).Value = 1957
$ws.Cells.Item(20,3).Value = 3069
$ws.Cells.Item(20,4).Value = 3007
$ws.Cells.Item(20,5).Value = 299
$ws.Cells.Item(20,6).Value = 5594
$ws.Cells.Item(20,7).Value = 2775
$ws.Cells.Item(20,8).Value = 179
$ws.Cells.Item(20,9).Value = 5316

$ws.Cells.Item(21,1).Value = "Sweden"
$ws.Cells.Item(21,2).Value = 1958
$ws.Cells.Item(21,3).Value = 1252
$ws.Cells.Item(21,4).Value = 783
$ws.Cells.Item(21,5).Value = -2333
$ws.Cells.Item(21,6).Value = 3741
$ws.Cells.Item(21,7).Value = -335
$ws.Cells.Item(21,8).Value = -3528
$ws.Cells.Item(21,9).Value = 2714

$ws.Cells.Item(22,1).Value = "Sweden"
$ws.Cells.Item(22,2).Value = 2019
$ws.Cells.Item(22,3).Value = -4114
$ws.Cells.Item(22,4).Value = -3960
$ws.Cells.Item(22,5).Value = -7751
$ws.Cells.Item(22,6).Value = -278
$ws.Cells.Item(22,7).Value = -4347
$ws.Cells.Item(22,8).Value = -8046
$ws.Cells.Item(22,9).Value = -850

$ws.Cells.Item(23,1).Value = "Sweden"
$ws.Cells.Item(23,2).Value = 2020
$ws.Cells.Item(23,3).Value = 8894
$ws.Cells.Item(23,4).Value = 8814
$ws.Cells.Item(23,5).Value = 5519
$ws.Cells.Item(23,6).Value = 12130
$ws.Cells.Item(23,7).Value = 7671
$ws.Cells.Item(23,8).Value = 4247
$ws.Cells.Item(23,9).Value = 11013

$ws.Cells.Item(24,1).Value = "Spain"
$ws.Cells.Item(24,2).Value = 1917
$ws.Cells.Item(24,3).Value = 12312
$ws.Cells.Item(24,4).Value = 11887
$ws.Cells.Item(24,5).Value = -12754
$ws.Cells.Item(24,6).Value = 36380
$ws.Cells.Item(24,7).Value = 11957
$ws.Cells.Item(24,8).Value = -13513
$ws.Cells.Item(24,9).Value = 35647

$ws.Cells.Item(25,1).Value = "Spain"
$ws.Cells.Item(25,2).Value = 1918
$ws.Cells.Item(25,3).Value = 243980
$ws.Cells.Item(25,4).Value = 244012
$ws.Cells.Item(25,5).Value = 218658
$ws.Cells.Item(25,6).Value = 267417
$ws.Cells.Item(25,7).Value = 242061
$ws.Cells.Item(25,8).Value = 216201
$ws.Cells.Item(25,9).Value = 267526

$ws.Cells.Item(26,1).Value = "Spain"
$ws.Cells.Item(26,2).Value = 1919
$ws.Cells.Item(26,3).Value = 26502
$ws.Cells.Item(26,4).Value = 25389
$ws.Cells.Item(26,5).Value = -2679
$ws.Cells.Item(26,6).Value = 51965
$ws.Cells.Item(26,7).Value = 28747
$ws.Cells.Item(26,8).Value = 1882
$ws.Cells.Item(26,9).Value = 54960

$ws.Cells.Item(27,1).Value = "Spain"
$ws.Cells.Item(27,2).Value = 1956
$ws.Cells.Item(27,3).Value = 39615
$ws.Cells.Item(27,4).Value = 33331
$ws.Cells.Item(27,5).Value = 10713
$ws.Cells.Item(27,6).Value = 54946
$ws.Cells.Item(27,7).Value = 32416
$ws.Cells.Item(27,8).Value = 9151
$ws.Cells.Item(27,9).Value = 53990

$ws.Cells.Item(28,1).Value = "Spain"
$ws.Cells.Item(28,2).Value = 1957
$ws.Cells.Item(28,3).Value = 10613
$ws.Cells.Item(28,4).Value = 11667
$ws.Cells.Item(28,5).Value = -5235
$ws.Cells.Item(28,6).Value = 27206
$ws.Cells.Item(28,7).Value = 13292
$ws.Cells.Item(28,8).Value = -3401
$ws.Cells.Item(28,9).Value = 28970

$ws.Cells.Item(29,1).Value = "Spain"
$ws.Cells.Item(29,2).Value = 1958
$ws.Cells.Item(29,3).Value = -33028
$ws.Cells.Item(29,4).Value = -29316
$ws.Cells.Item(29,5).Value = -49105
$ws.Cells.Item(29,6).Value = -10269
$ws.Cells.Item(29,7).Value = -28860
$ws.Cells.Item(29,8).Value = -48158
$ws.Cells.Item(29,9).Value = -9413

$ws.Cells.Item(30,1).Value = "Spain"
$ws.Cells.Item(30,2).Value = 2019
$ws.Cells.Item(30,3).Value = -21233
$ws.Cells.Item(30,4).Value = -20616
$ws.Cells.Item(30,5).Value = -42795
$ws.Cells.Item(30,6).Value = 761
$ws.Cells.Item(30,7).Value = -13200
$ws.Cells.Item(30,8).Value = -33870
$ws.Cells.Item(30,9).Value = 7717

$ws.Cells.Item(31,1).Value = "Spain"
$ws.Cells.Item(31,2).Value = 2020
$ws.Cells.Item(31,3).Value = 67428
$ws.Cells.Item(31,4).Value = 66826
$ws.Cells.Item(31,5).Value = 45817
$ws.Cells.Item(31,6).Value = 87210
$ws.Cells.Item(31,7).Value = 70926
$ws.Cells.Item(31,8).Value = 50451
$ws.Cells.Item(31,9).Value = 90911
